# Update "想去人数" (want-to-go count) figures across sheets after a data refresh.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# 展览 (sheet1)
$ws1.Range("F2").Value = 5714
$ws1.Range("F3").Value = 7684
$ws1.Range("F4").Value = 82
$ws1.Range("F6").Value = 616
$ws1.Range("F9").Value = 4483
$ws1.Range("F12").Value = 123
$ws1.Range("F13").Value = 3023
$ws1.Range("F16").Value = 229
$ws1.Range("F17").Value = 568
$ws1.Range("F19").Value = 490
$ws1.Range("F22").Value = 1739
$ws1.Range("F23").Value = 1271
$ws1.Range("F24").Value = 108
$ws1.Range("F25").Value = 1482
$ws1.Range("F31").Value = 30
$ws1.Range("F32").Value = 76
$ws1.Range("F34").Value = 78
$ws1.Range("F35").Value = 3265
$ws1.Range("F37").Value = 51
$ws1.Range("F38").Value = 200
$ws1.Range("F40").Value = 1325

# 演出 (sheet2)
$ws2.Range("F2").Value = 17
$ws2.Range("F3").Value = 29

# 全部类型 (sheet4) - aggregated view, updated with same figures
$ws4.Range("F2").Value = 5714
$ws4.Range("F3").Value = 7684
$ws4.Range("F4").Value = 82
$ws4.Range("F6").Value = 616
$ws4.Range("F9").Value = 4483
$ws4.Range("F12").Value = 123
$ws4.Range("F13").Value = 3023
$ws4.Range("F16").Value = 229
$ws4.Range("F17").Value = 568
$ws4.Range("F19").Value = 490
$ws4.Range("F20").Value = 17
$ws4.Range("F23").Value = 1739
$ws4.Range("F24").Value = 1271
$ws4.Range("F25").Value = 108
$ws4.Range("F26").Value = 1482
$ws4.Range("F32").Value = 30
$ws4.Range("F33").Value = 76
$ws4.Range("F35").Value = 78
$ws4.Range("F36").Value = 3265
$ws4.Range("F37").Value = 29
$ws4.Range("F39").Value = 51
$ws4.Range("F40").Value = 200
$ws4.Range("F42").Value = 1325
